# Applies the "First py ta fm" commit:
#  - B2: "TC001 SearchInJioMart" -> "TC001 TestCase_001"
#  - D2: "Search Item in Jio mart " -> "Test Case 1 "
#  - A2 formula recalculates to "TC_TC001_TestCase_001" automatically
#  - active selection moves to C2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "TC001 TestCase_001"
$ws.Range("D2").Value = "Test Case 1 "

$ws.Range("C2").Select()
